$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rhamnMat_corrected")

# Relabel strain W36662 (row 26) as a non-producer: set rhamn3cats and rhamn2cats to 0
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0

# Mirror the selection left behind after the edit (cell C28 active)
$ws.Range("C28").Select()
